# contratos-11-2013.xlsx fix (per commit "fix: fixed formatting when
# scrapping floating point numbers"):
#
#  1) A handful of "Razon social"/"Nombre Fantasia" text entries (column
#     E / F) had a comma used where the scraper should have used a
#     period (e.g. "SCHAB DARIO, PEROTTI XAVIER, BENINCA MATIAS S.H."
#     -> "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"). Fix the typo
#     in every cell that carries the bad text.
#
#  2) Every "Importe" value in column H was scraped as European-style
#     text ("1.040.300,00" - "." thousands separator, "," decimal
#     separator). The fix re-writes each one as plain-decimal text
#     ("1040300.00" - no thousands separator, "." decimal separator).
#     The cells stay plain text (not real numbers), so a leading
#     apostrophe is used to stop Excel's auto-conversion from turning
#     the digits back into a Number when the value is assigned.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 1) Razon social / Nombre Fantasia comma -> period fixes ----------
$textFixes = @{
    "TRABICHET MARIA. VERGARA ADEL Y OTRA" = @("E31","F31","E125","F125","E129","F129")
    "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA" = @("E38","F38","E55","F55","E121","F121","E185","F185")
    "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO" = @("E93","E179")
    "RICCOTTI. MARIANA EDITH" = @("E137")
    "MERCANZINI. GASTON ARIEL" = @("F149")
    "DODERA. JORGE ABELARDO" = @("E183")
    "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH" = @("E188","E209","E237","E244")
    "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN" = @("E191")
}

foreach ($newText in $textFixes.Keys) {
    foreach ($addr in $textFixes[$newText]) {
        $ws.Range($addr).Value = $newText
    }
}

# ---- 2) Importe (column H) European number-text -> plain number-text --
$firstRow = 2
$lastRow = 253
$col = 8  # column H

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $old = $cell.Value2
    if ($old -eq $null) { continue }
    $new = $old.Replace(".", "").Replace(",", ".")
    # Leading apostrophe keeps this a text value (matches original t="s"
    # shared-string cell) instead of letting Excel re-parse it as a Number.
    $cell.Value = "'" + $new
}
